$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")
$ws.Range("A12").Value = "Source"
$ws.Range("B12").Value = "NATSISS 2008 and 2014/15"
$ws.Range("B13").Value = "State and Territory governments"
$ws.Range("A13").Select() | Out-Null
